$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "answers"
$ws.Range("D1").Value = "topic"
$ws.Range("E1").Value = "subtopic"
$ws.Range("F1").Value = "difficulty"
$ws.Range("A1").Value = "questionId"
$ws.Range("B1").Value = "prompt"

$ws.Range("F1").Select() | Out-Null
